# Update the "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1076
    3  = 361
    4  = 1473
    5  = 8695
    6  = 87
    7  = 489
    8  = 639
    11 = 9
    12 = 3542
    16 = 1139
    17 = 145
    18 = 1111
    21 = 2272
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
